# Revise controlled vocabulary public-services-subject-matters (issue #177)
# Merge/rename a few top-level subject-matter categories:
#   - "Agricoltura" -> "Agricoltura e pesca"
#   - "Tributi e finanze" -> "Tributi, finanze e contravvenzioni" (and tweak its long label)
#   - "Attività produttive e commercio" -> "Impresa e commercio" (and tweak its long label)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - Agricoltura -> Agricoltura e pesca
$ws.Range("B16").Value = "Agricoltura e pesca"

# Row 10 - Tributi e finanze -> Tributi, finanze e contravvenzioni
$ws.Range("B10").Value = "Tributi, finanze e contravvenzioni"

# Row 13 - Attività produttive e commercio -> Impresa e commercio
$ws.Range("B13").Value = "Impresa e commercio"
$ws.Range("C13").Value = "Impresa e commercio - attività produttive, impresa nazionale/estera, notifiche, bancarotta, risorse umane"

# Row 16 alternative label
$ws.Range("C16").Value = "Agricoltura e pesca - politiche agricole e alimentari"

# Row 10 alternative label
$ws.Range("C10").Value = "Tributi, finanze e contravvenzioni -  dichiarazione redditi, contributi"

$ws.Range("C12").Select()
